$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 594, shifting existing rows 594-658 down to 595-659
$ws.Rows.Item(594).Insert()

# Populate the newly inserted row 594 with the new record's data.
# Columns that stay constant across every record in this sheet (A,B,C,E,F,G,H,I,J,K,L,Q,R,T)
# are carried over unchanged; D,M,N,O,P,S are the new values for this entry.
$ws.Cells.Item(594, 1).Value  = 4
$ws.Cells.Item(594, 2).Value  = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(594, 3).Value  = "Los Lagos"
$ws.Cells.Item(594, 4).Value  = 45142
$ws.Cells.Item(594, 5).Value  = 10
$ws.Cells.Item(594, 6).Value  = "Fruta"
$ws.Cells.Item(594, 7).Value  = 100102
$ws.Cells.Item(594, 8).Value  = "Cítricos"
$ws.Cells.Item(594, 9).Value  = 100102006
$ws.Cells.Item(594, 10).Value = "Pomelo"
$ws.Cells.Item(594, 11).Value = "Start Ruby"
$ws.Cells.Item(594, 12).Value = "Primera"
$ws.Cells.Item(594, 13).Value = 150
$ws.Cells.Item(594, 14).Value = 14000
$ws.Cells.Item(594, 15).Value = 14000
$ws.Cells.Item(594, 16).Value = 14000
$ws.Cells.Item(594, 17).Value = '$/caja 14 kilos empedrada'
$ws.Cells.Item(594, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(594, 19).Value = 1000
$ws.Cells.Item(594, 20).Value = 14
